$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current "jx:if" row (row 5), which pushes
# the jx:if row down to row 6 and the "Parent" row down to row 7.
$ws.Rows.Item(5).Insert()

# Split "Name: ${person.name}" into a label cell + a value cell.
$ws.Range("A3").Value = "Name:"
$ws.Range("B3").Value = '${person.name}'

# Split "Age: ${person.age}" into a label cell + a value cell.
$ws.Range("A4").Value = "Age:"
$ws.Range("B4").Value = '${person.age}'

# Update the jx:if command row (now row 6) with the new condition syntax.
$ws.Range("A6").Value = 'jx:if(condition="person.age < 18", lastCell="B6")'

# Split "Parent: ${person.parentName}" into a label cell + a value cell
# (now on row 7 after the row insert).
$ws.Range("A7").Value = "Parent:"
$ws.Range("B7").Value = '${person.parentName}'

# Resize columns: A narrower, new column B added for values.
$ws.Columns.Item(1).ColumnWidth = 14.791666666666666
$ws.Columns.Item(2).ColumnWidth = 22.604166666666668

# Attach an explanatory comment to the jx:if cell.
$comment = $ws.Range("A6").AddComment("JXLS Command:`nThis row contains the jx:if condition.`nIf condition is true, the next row(s) will be included.`nIf false, they will be removed from output.")

# Recolor the jx:if row's fill (indexed 22 -> indexed 43, a pale yellow).
$ws.Range("A6").Interior.Color = 10092543
